# Generate Report for Handoff
# Replaces the two tracked e2e markdown docs (old GUIDs) with a new pair of
# GUIDs, flips the localization status from "Handed back: in sync with
# en-US" to "Ready for handoff", refreshes the handoff timestamps / xliff
# file names, and clears the now-stale "Latest Target File" columns.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "7678ce99-c951-4420-978b-50240a635986"
$oldGuid2 = "82bf6c72-145b-4e01-8ed8-41f688a28fec"
$newGuid1 = "05049e00-37bb-4c68-ae9a-126150ae4e7f"
$newGuid2 = "ffff0c29e0a4-cc9d-4261-bd52-5c0058008b0d"

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newGuid1.md"
$wsOverview.Range("B2").Value2 = "e2e\$newGuid1.md"
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("G2").Value2 = "2016-08-15 14:59:51"

$wsOverview.Range("A3").Value2 = "$newGuid2.md"
$wsOverview.Range("B3").Value2 = "e2e\$newGuid2.md"
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus
$wsOverview.Range("G3").Value2 = "2016-08-15 14:59:51"

# Rebuild the two row hyperlinks (B2 / B3) with the new guids.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

# Columns E/F got narrower now that the status text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = "$newGuid1.md"
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("G2").Value2 = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$wsZhCn.Range("H2").Value2 = "2016-08-15 14:59:44"
$wsZhCn.Range("I2").Value2 = ""
$wsZhCn.Range("J2").Value2 = ""
$wsZhCn.Range("K2").Value2 = "0001-01-01 00:00:00"

$wsZhCn.Range("A3").Value2 = "$newGuid2.md"
$wsZhCn.Range("C3").Value2 = $newStatus
$wsZhCn.Range("F3").Value2 = "True"
$wsZhCn.Range("G3").Value2 = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.zh-cn.xlf"
$wsZhCn.Range("H3").Value2 = "2016-08-15 14:59:44"
$wsZhCn.Range("I3").Value2 = ""
$wsZhCn.Range("J3").Value2 = ""
$wsZhCn.Range("K3").Value2 = "0001-01-01 00:00:00"

# The "Latest Target File" hyperlinks (I2/I3) are gone now that those cells
# are blank; only the Source File Name (A2/A3) links remain.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I3").Style = "Normal"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.38
$wsZhCn.Columns.Item(9).ColumnWidth = 17.83
$wsZhCn.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = "$newGuid1.md"
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("G2").Value2 = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$wsDeDe.Range("H2").Value2 = "2016-08-15 14:59:51"
$wsDeDe.Range("I2").Value2 = ""
$wsDeDe.Range("J2").Value2 = ""
$wsDeDe.Range("K2").Value2 = "0001-01-01 00:00:00"

$wsDeDe.Range("A3").Value2 = "$newGuid2.md"
$wsDeDe.Range("C3").Value2 = $newStatus
$wsDeDe.Range("F3").Value2 = "True"
$wsDeDe.Range("G3").Value2 = "$newGuid1.62faf983f671d06f1b12dade6646264544e38dfe.de-de.xlf"
$wsDeDe.Range("H3").Value2 = "2016-08-15 14:59:51"
$wsDeDe.Range("I3").Value2 = ""
$wsDeDe.Range("J3").Value2 = ""
$wsDeDe.Range("K3").Value2 = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d480cdede9972998ba2e5b3899569272dbb94fa/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I3").Style = "Normal"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.38
$wsDeDe.Columns.Item(9).ColumnWidth = 17.83
$wsDeDe.Columns.Item(10).ColumnWidth = 20.83
